$wb = $excel.ActiveWorkbook

# Update "展览" sheet (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 25
$ws1.Range("F4").Value = 982

# Update "全部类型" sheet (all types) - mirrors the same data
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 25
$ws4.Range("F4").Value = 982
